$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 51 values (set in the same order the original authors' shared
# strings appear in the diff: OPQA-2118, Profile50, Verify...)
$ws.Range("B51").Value = "OPQA-2118"
$ws.Range("A51").Value = "Profile50"
$ws.Range("C51").Value = "Verify that system is capturing the on-boarded events."
$ws.Range("D51").Value = "Y"

# Copy existing cell formatting (border/fill/alignment/font) onto the new
# row so styles are reused rather than re-invented.
$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial(-4122)

$ws.Range("B49").Copy()
$ws.Range("B51").PasteSpecial(-4122)

$ws.Range("C50").Copy()
$ws.Range("C51").PasteSpecial(-4122)

$ws.Range("D49").Copy()
$ws.Range("D51").PasteSpecial(-4122)

$ws.Range("E50").Copy()
$ws.Range("E51").PasteSpecial(-4122)

# C51 uses a distinct font (Arial, black) compared to C50.
$ws.Range("C51").Font.Name = "Arial"
$ws.Range("C51").Font.ColorIndex = 1

$excel.CutCopyMode = 0

# Update the active selection recorded in the sheet view, as Excel does
# when the user clicks a different cell before saving.
[void]$ws.Range("C43").Select()
